# Update the Top 10 Fortune Global 500 (2019) data and its panel view:
#   - "State Grid" is renamed to "State Grid (China)" to match the clarified
#     company naming used alongside the accompanying maps.
#   - The worksheet's active selection/window view is moved to reflect the
#     author's latest working position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the company in row 6 (rank 5) from "State Grid" to "State Grid (China)"
$ws.Range("B6").Value = "State Grid (China)"

# Reposition the workbook window (best effort - matches author's window move)
$win = $wb.Windows.Item(1)
$win.Left = 1560
$win.Top = 0

# Update the active selection on the worksheet to B7
$ws.Activate()
[void]($ws.Range("B7").Select())
